# Refresh Leve profit-calc columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# per updated market-board pricing data -- scheduled runner sync.
$wb = $excel.ActiveWorkbook


# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item(1)
# Row 6: Days of Chunder / Antidote
$ws.Range("H6").Value = 1111198.8
$ws.Range("I6").Value = 1428645.6
$ws.Range("K6").Value = 4285936.800000001
$ws.Range("M6").Value = -4285824.800000001
# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 3930
$ws.Range("I12").Value = 6350
$ws.Range("K12").Value = 6350
$ws.Range("M12").Value = -6180
# Row 32: Automata for the People / Crab Oil
$ws.Range("H32").Value = 8668.75
$ws.Range("J32").Value = 9891.666999999999
$ws.Range("L32").Value = 9891.666999999999
$ws.Range("N32").Value = -10543.667
# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 8000836.5
$ws.Range("I92").Value = 8000836.5
$ws.Range("K92").Value = 8000836.5
$ws.Range("M92").Value = -7999588.5
# Row 97: Materia Worth / Potent Spiritbond Potion
$ws.Range("H97").Value = 58372
$ws.Range("J97").Value = 166666
$ws.Range("L97").Value = 499998
$ws.Range("N97").Value = -500990
# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 1121.7858
$ws.Range("I100").Value = 1148.2
$ws.Range("K100").Value = 1148.2
$ws.Range("M100").Value = -607.2
# Row 104: Pep-stepper / Infusion of Vitality
$ws.Range("H104").Value = 250.16667
$ws.Range("I104").Value = 265.75
$ws.Range("J104").Value = 219
$ws.Range("K104").Value = 797.25
$ws.Range("L104").Value = 657
$ws.Range("M104").Value = 949.75
$ws.Range("N104").Value = -4151
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 6538377
$ws.Range("I132").Value = 7937900.5
$ws.Range("J132").Value = 7266.5557
$ws.Range("K132").Value = 23813701.5
$ws.Range("L132").Value = 21799.6671
$ws.Range("M132").Value = -23811171.5
$ws.Range("N132").Value = -26859.6671
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2695.9033
$ws.Range("I137").Value = 2386.5
$ws.Range("K137").Value = 7159.5
$ws.Range("M137").Value = -4609.5
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 7694711
$ws.Range("I138").Value = 962.2222
$ws.Range("J138").Value = 13161322
$ws.Range("K138").Value = 2886.6666
$ws.Range("L138").Value = 39483966
$ws.Range("M138").Value = 2253.3334
$ws.Range("N138").Value = -39494246
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 1245.1333
$ws.Range("I141").Value = 1245.1333
$ws.Range("K141").Value = 3735.3999
$ws.Range("M141").Value = 1444.6001

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item(2)
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 4373.9316
$ws.Range("I32").Value = 2153.8955
$ws.Range("K32").Value = 2153.8955
$ws.Range("M32").Value = -1866.8955
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2479.986
$ws.Range("I61").Value = 1516.5555
$ws.Range("J61").Value = 5370.278
$ws.Range("K61").Value = 1516.5555
$ws.Range("L61").Value = 5370.278
$ws.Range("M61").Value = -1304.5555
$ws.Range("N61").Value = -5794.278
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 27968.822
$ws.Range("I74").Value = 28767.703
$ws.Range("J74").Value = 24274
$ws.Range("K74").Value = 28767.703
$ws.Range("L74").Value = 24274
$ws.Range("M74").Value = -27893.703
$ws.Range("N74").Value = -26022
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 27968.822
$ws.Range("I77").Value = 28767.703
$ws.Range("J77").Value = 24274
$ws.Range("K77").Value = 143838.515
$ws.Range("L77").Value = 121370
$ws.Range("M77").Value = -139470.515
$ws.Range("N77").Value = -130106
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 1860.3636
$ws.Range("I97").Value = 1501.2354
$ws.Range("K97").Value = 1501.2354
$ws.Range("M97").Value = -1005.2354
# Row 123: The Armoire Is Open / High Durium Armguards of Maiming
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2677.4312
$ws.Range("I132").Value = 2538.2
$ws.Range("J132").Value = 3159.3845
$ws.Range("K132").Value = 7614.599999999999
$ws.Range("L132").Value = 9478.1535
$ws.Range("M132").Value = -5084.599999999999
$ws.Range("N132").Value = -14538.1535
# Row 134: Brace for More Vambraces / Ruthenium Vambraces of Maiming
$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2479.986
$ws.Range("I136").Value = 1516.5555
$ws.Range("J136").Value = 5370.278
$ws.Range("K136").Value = 4549.666499999999
$ws.Range("L136").Value = 16110.834
$ws.Range("M136").Value = -1999.666499999999
$ws.Range("N136").Value = -21210.834

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item(3)
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2382.6667
$ws.Range("I105").Value = 2474
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 2474
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -727
$ws.Range("N105").Value = -5694

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item(4)
# Row 3: Touch and Heal / Maple Pattens
$ws.Range("H3").Value = 24120
$ws.Range("I3").Value = 5150
$ws.Range("K3").Value = 5150
$ws.Range("M3").Value = -5037
# Row 4: A Clogful of Camaraderie / Maple Clogs
$ws.Range("H4").Value = 1066085.8
$ws.Range("I4").Value = 1389520
$ws.Range("J4").Value = 257500
$ws.Range("K4").Value = 1389520
$ws.Range("L4").Value = 257500
$ws.Range("M4").Value = -1389408
$ws.Range("N4").Value = -257724
# Row 20: Re-crating the Scene / Iron Spear
$ws.Range("H20").Value = 250000
$ws.Range("J20").Value = 250000
$ws.Range("L20").Value = 250000
$ws.Range("N20").Value = -250472
# Row 30: Polearms Aplenty / Iron Spear
$ws.Range("H30").Value = 250000
$ws.Range("J30").Value = 250000
$ws.Range("L30").Value = 250000
$ws.Range("N30").Value = -250182
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 18911.8
$ws.Range("I31").Value = 26576.285
$ws.Range("J31").Value = 4915.7827
$ws.Range("K31").Value = 26576.285
$ws.Range("L31").Value = 4915.7827
$ws.Range("M31").Value = -26281.285
$ws.Range("N31").Value = -5505.7827
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 18911.8
$ws.Range("I34").Value = 26576.285
$ws.Range("J34").Value = 4915.7827
$ws.Range("K34").Value = 26576.285
$ws.Range("L34").Value = 4915.7827
$ws.Range("M34").Value = -26374.285
$ws.Range("N34").Value = -5319.7827
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 2646.843
$ws.Range("I58").Value = 2417.413
$ws.Range("J58").Value = 4757.6
$ws.Range("K58").Value = 2417.413
$ws.Range("L58").Value = 4757.6
$ws.Range("M58").Value = -2214.413
$ws.Range("N58").Value = -5163.6
# Row 60: Bowing to Greater Power / Yew Longbow
$ws.Range("H60").Value = 16874.25
$ws.Range("J60").Value = 19999
$ws.Range("L60").Value = 19999
$ws.Range("N60").Value = -21021
# Row 110: A Stronger Offense / Applewood Spear
$ws.Range("H110").Value = 250000
$ws.Range("J110").Value = 250000
$ws.Range("L110").Value = 250000
$ws.Range("N110").Value = -258180
# Row 115: Horde of the Rings / White Ash Bracelet
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
# Row 116: The Right Tool for the Job / Sandteak Rod
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 128: An A-prop-riate Request / Ironwood Spear
$ws.Range("H128").Value = 250000
$ws.Range("J128").Value = 250000
$ws.Range("L128").Value = 250000
$ws.Range("N128").Value = -259960
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 2646.843
$ws.Range("I136").Value = 2417.413
$ws.Range("J136").Value = 4757.6
$ws.Range("K136").Value = 7252.239
$ws.Range("L136").Value = 14272.8
$ws.Range("M136").Value = -4702.239
$ws.Range("N136").Value = -19372.8

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item(5)
# Row 6: Meat-lover's Special / Marmot Steak
$ws.Range("H6").Value = 353.41666
$ws.Range("I6").Value = 353.41666
$ws.Range("K6").Value = 1060.24998
$ws.Range("M6").Value = -947.2499800000001
# Row 10: A Real Fungi / Chanterelle Saute
$ws.Range("H10").Value = 37.6
$ws.Range("I10").Value = 37.6
$ws.Range("K10").Value = 112.8
$ws.Range("M10").Value = 26.19999999999999
# Row 43: Sole Survivor / Baked Sole
$ws.Range("H43").Value = 5000.6665
$ws.Range("J43").Value = 5500
$ws.Range("L43").Value = 16500
$ws.Range("N43").Value = -16728
# Row 55: Pagan Pastries / Pastry Fish
$ws.Range("H55").Value = 2823.111
$ws.Range("J55").Value = 3492.8572
$ws.Range("L55").Value = 10478.5716
$ws.Range("N55").Value = -10832.5716
# Row 117: A Good Omen / Peppered Popotoes
$ws.Range("H117").Value = 670.9375
$ws.Range("I117").Value = 1014.6
$ws.Range("K117").Value = 3043.8
$ws.Range("M117").Value = 398.1999999999998

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item(6)
# Row 3: Needful Rings / Copper Wristlets
$ws.Range("H3").Value = 25001946
$ws.Range("I3").Value = 33335426
$ws.Range("K3").Value = 33335426
$ws.Range("M3").Value = -33335310
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2472.0178
$ws.Range("I132").Value = 2449.6738
$ws.Range("J132").Value = 2574.8
$ws.Range("K132").Value = 7349.0214
$ws.Range("L132").Value = 7724.400000000001
$ws.Range("M132").Value = -4819.0214
$ws.Range("N132").Value = -12784.4

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item(7)
# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 19833
$ws.Range("I100").Value = 7000
$ws.Range("K100").Value = 7000
$ws.Range("M100").Value = -6459
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4016.1
$ws.Range("I132").Value = 4100
$ws.Range("K132").Value = 12300
$ws.Range("M132").Value = -9770
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4250.276
$ws.Range("I136").Value = 3550.6191
$ws.Range("J136").Value = 6086.875
$ws.Range("K136").Value = 10651.8573
$ws.Range("L136").Value = 18260.625
$ws.Range("M136").Value = -8101.8573
$ws.Range("N136").Value = -23360.625

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item(8)
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 8981.368
$ws.Range("I81").Value = 21659.4
$ws.Range("J81").Value = 4453.5
$ws.Range("K81").Value = 43318.8
$ws.Range("L81").Value = 8907
$ws.Range("M81").Value = -42257.8
$ws.Range("N81").Value = -11029
# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 8981.368
$ws.Range("I84").Value = 21659.4
$ws.Range("J84").Value = 4453.5
$ws.Range("K84").Value = 216594
$ws.Range("L84").Value = 44535
$ws.Range("M84").Value = -211290
$ws.Range("N84").Value = -55143
# Row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 1808.2
$ws.Range("I96").Value = 1548.2307
$ws.Range("K96").Value = 1548.2307
$ws.Range("M96").Value = -175.2307000000001
